$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Table 2 ("Operation Name" / "Resource URL " / ... ), first data row,
#    first cell: "[a.api_name] " is currently a single run. Split it into
#    three runs: "[a." + "op" + "name] " (same character formatting on all
#    three) -- this mirrors what Word's spell-checker/autocorrect does when
#    it silently fragments a run while you are typing/correcting text.
# ---------------------------------------------------------------------------
$tbl = $d.Tables.Item(2)
$cell = $tbl.Rows.Item(2).Cells.Item(1)
$cellStart = $cell.Range.Start

# Original text layout inside the cell, relative to $cellStart:
#   [0,3)  -> "[a."
#   [3,5)  -> "ap"
#   [5,13) -> "i_name] "
$seg1 = $d.Range($cellStart, $cellStart + 3)
$null = $seg1.Find.Execute("[a.", $false, $false, $false, $false, $false, $true, 1, $false, "[a.", 2)

# Toggling Bold on/off around the replace forces Word to materialize a new,
# independent run instead of just editing the text of the existing run in
# place -- then resetting Bold back to its original value removes the
# formatting difference again, leaving a clean run boundary behind.
$seg2 = $d.Range($cellStart + 3, $cellStart + 5)
$seg2.Bold = 1
$null = $seg2.Find.Execute("ap", $false, $false, $false, $false, $false, $true, 1, $false, "op", 2)
$seg2b = $d.Range($cellStart + 3, $cellStart + 5)
$seg2b.Bold = 0

$seg3 = $d.Range($cellStart + 5, $cellStart + 13)
$seg3.Bold = 1
$null = $seg3.Find.Execute("i_name] ", $false, $false, $false, $false, $false, $true, 1, $false, "name] ", 2)
$seg3b = $d.Range($cellStart + 5, $cellStart + 11)
$seg3b.Bold = 0

# ---------------------------------------------------------------------------
# 2) Remove the "OpenTBS template tips" block together with the
#    "Merging data with a Chart" heading and its chart illustration, right
#    before the "Merging data with a table" section -- keeping the single
#    blank paragraph that precedes them.
# ---------------------------------------------------------------------------
$rng = $d.Content
$null = $rng.Find.Execute("You may consider this before building your own Microsoft Word template:", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$deleteStart = $rng.Start

$rng2 = $d.Content
$null = $rng2.Find.Execute("Merging data with a table", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$deleteEnd = $rng2.Start

$d.Range($deleteStart, $deleteEnd).Delete()
